# Normalize the "Updated By" (column G) entries on the "Session Analysis
# Results" sheet: for any cell whose value is exactly two comma-separated
# names/emails, put them in the opposite order (e.g. "System" moves to the
# front when it was second). Cells with a single entry, or with three or
# more entries, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ','
    if ($parts.Count -eq 2) {
        $first = $parts[0].Trim()
        $second = $parts[1].Trim()

        if ($first -ne 'System') {
            $cell.Value = "$second, $first"
        }
    }
}
